# Update entregable 1, 2 y 3
# Refresh the "control dual" sample row with the latest client data:
#   - mnemocino / dni now reference client 11122548 (was 75356819)
#   - nacimiento (birth date) corrected to 19980101 (was 19980112)
#   - Fecha (approval timestamp) advanced to 2 jul. 2023, 17:15:22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# mnemocino (text, starts with a letter so it stores naturally as text)
$ws.Range("C2").Value = "D11122548"

# dni / nacimiento look like numbers -- use a leading apostrophe (quote-prefix)
# so they are kept as text, matching how they were originally authored.
$ws.Range("D2").Formula = "'11122548"
$ws.Range("J2").Formula = "'19980101"

# Fecha (plain text timestamp)
$ws.Range("P2").Value = "2 jul. 2023, 17:15:22"

# The dni column (D) now holds data, so best-fit its width like the
# other data columns on the sheet.
$ws.Columns.Item(4).AutoFit()

# Leave the cursor where the edit was made.
$null = $ws.Range("D9").Select()
